# 自动更新价格数据: insert a new "latest" row (2025-12-16) at the top of the
# data table (row 2), pushing the existing history down by one row and
# keeping the previous bottom row (2025-11-21) in the table as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last used row of the current table (row 26 in the starting workbook: header
# in row 1 + 25 data rows).
$lastRow = $ws.UsedRange.Rows.Count

# Column A holds dates stored as plain text (e.g. "2025-12-15"), not real
# Excel dates. Force text formatting on the rows we are about to (re)write so
# the "yyyy-mm-dd" strings aren't auto-converted into date serial numbers.
$textRange = "A2:A$($lastRow + 1)"
$ws.Range($textRange).NumberFormat = "@"

# Shift every existing data row down by one (start from the bottom so we
# don't clobber a row before it has been copied).
for ($r = $lastRow; $r -ge 2; $r--) {
  $dst = $r + 1
  $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
  $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
  $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value2
  $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value2
}

# Write the new latest-day row.
$ws.Cells.Item(2, 1).Value = "2025-12-16"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
